$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 29412996
$ws.Range("I43").Value = 83334330
$ws.Range("J43").Value = 1361.909
$ws.Range("K43").Value = 83334330
$ws.Range("L43").Value = 1361.909
$ws.Range("M43").Value = -83334261
$ws.Range("N43").Value = -1499.909

$ws.Range("H109").Value = 13762
$ws.Range("J109").Value = 13762
$ws.Range("L109").Value = 13762
$ws.Range("N109").Value = -16536

$ws.Range("H115").Value = 263.33334
$ws.Range("I115").Value = 263.33334
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 790.0000200000001
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 776.9999799999999
$ws.Range("N115").ClearContents()

$ws.Range("H129").Value = 943.7681
$ws.Range("J129").Value = 912.0476
$ws.Range("L129").Value = 2736.1428
$ws.Range("N129").Value = -12736.1428

$ws.Range("H132").Value = 18756920
$ws.Range("I132").Value = 21436138
$ws.Range("J132").Value = 2399.2
$ws.Range("K132").Value = 64308414
$ws.Range("L132").Value = 7197.599999999999
$ws.Range("M132").Value = -64305884
$ws.Range("N132").Value = -12257.6

$ws.Range("H138").Value = 5840.6724
$ws.Range("I138").Value = 1659.2858
$ws.Range("J138").Value = 9743.299999999999
$ws.Range("K138").Value = 4977.857400000001
$ws.Range("L138").Value = 29229.9
$ws.Range("M138").Value = 162.1425999999992
$ws.Range("N138").Value = -39509.89999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 3000
$ws.Range("I34").Value = 3000
$ws.Range("K34").Value = 3000
$ws.Range("M34").Value = -2729

$ws.Range("H54").Value = 22000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 22000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 22000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -23538

$ws.Range("H88").Value = 3601.7144
$ws.Range("I88").Value = 4240.6665
$ws.Range("J88").Value = 3122.5
$ws.Range("K88").Value = 4240.6665
$ws.Range("L88").Value = 3122.5
$ws.Range("M88").Value = -3834.6665
$ws.Range("N88").Value = -3934.5

$ws.Range("H91").Value = 3601.7144
$ws.Range("I91").Value = 4240.6665
$ws.Range("J91").Value = 3122.5
$ws.Range("K91").Value = 4240.6665
$ws.Range("L91").Value = 3122.5
$ws.Range("M91").Value = -2836.6665
$ws.Range("N91").Value = -5930.5

$ws.Range("H121").Value = 24399.5
$ws.Range("J121").Value = 24399.5
$ws.Range("L121").Value = 24399.5
$ws.Range("N121").Value = -27893.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 49800
$ws.Range("J59").Value = 49800
$ws.Range("L59").Value = 49800
$ws.Range("N59").Value = -51494

$ws.Range("H86").Value = 1369356.8
$ws.Range("I86").Value = 947
$ws.Range("J86").Value = 1790405.9
$ws.Range("K86").Value = 947
$ws.Range("L86").Value = 1790405.9
$ws.Range("M86").Value = 176
$ws.Range("N86").Value = -1792651.9

$ws.Range("H89").Value = 1369356.8
$ws.Range("I89").Value = 947
$ws.Range("J89").Value = 1790405.9
$ws.Range("K89").Value = 4735
$ws.Range("L89").Value = 8952029.5
$ws.Range("M89").Value = 881
$ws.Range("N89").Value = -8963261.5

$ws.Range("H94").Value = 1000
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 1000
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -549
$ws.Range("N94").Value = -1902

$ws.Range("H134").Value = 1953289.4
$ws.Range("I134").Value = 1626.5116
$ws.Range("J134").Value = 7947682
$ws.Range("K134").Value = 4879.5348
$ws.Range("L134").Value = 23843046
$ws.Range("M134").Value = -2344.5348
$ws.Range("N134").Value = -23848116

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3624233.2
$ws.Range("I31").Value = 761.7714
$ws.Range("J31").Value = 15153461
$ws.Range("K31").Value = 761.7714
$ws.Range("L31").Value = 15153461
$ws.Range("M31").Value = -466.7714
$ws.Range("N31").Value = -15154051

$ws.Range("H34").Value = 3624233.2
$ws.Range("I34").Value = 761.7714
$ws.Range("J34").Value = 15153461
$ws.Range("K34").Value = 761.7714
$ws.Range("L34").Value = 15153461
$ws.Range("M34").Value = -559.7714
$ws.Range("N34").Value = -15153865

$ws.Range("H122").Value = 13159476
$ws.Range("I122").Value = 15626447
$ws.Range("J122").Value = 2298
$ws.Range("K122").Value = 46879341
$ws.Range("L122").Value = 6894
$ws.Range("M122").Value = -46876891
$ws.Range("N122").Value = -11794

$ws.Range("H132").Value = 9525441
$ws.Range("I132").Value = 1121.6296
$ws.Range("J132").Value = 41670020
$ws.Range("K132").Value = 3364.8888
$ws.Range("L132").Value = 125010060
$ws.Range("M132").Value = -834.8887999999997
$ws.Range("N132").Value = -125015120

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2783.5
$ws.Range("I68").Value = 710.6667
$ws.Range("J68").Value = 3880.8823
$ws.Range("K68").Value = 2132.0001
$ws.Range("L68").Value = 11642.6469
$ws.Range("M68").Value = -1321.0001
$ws.Range("N68").Value = -13264.6469

$ws.Range("H71").Value = 2783.5
$ws.Range("I71").Value = 710.6667
$ws.Range("J71").Value = 3880.8823
$ws.Range("K71").Value = 6396.0003
$ws.Range("L71").Value = 34927.9407
$ws.Range("M71").Value = -2340.0003
$ws.Range("N71").Value = -43039.9407

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4250
$ws.Range("I126").Value = 3100
$ws.Range("J126").Value = 4742.857
$ws.Range("K126").Value = 9300
$ws.Range("L126").Value = 14228.571
$ws.Range("M126").Value = -6830
$ws.Range("N126").Value = -19168.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1500
$ws.Range("I7").Value = 1500
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1388
$ws.Range("N7").ClearContents()

$ws.Range("H126").Value = 1500
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2030
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 16397284
$ws.Range("I132").Value = 24391776
$ws.Range("J132").Value = 8575.4
$ws.Range("K132").Value = 73175328
$ws.Range("L132").Value = 25726.2
$ws.Range("M132").Value = -73172798
$ws.Range("N132").Value = -30786.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 18000
$ws.Range("J51").Value = 18000
$ws.Range("L51").Value = 18000
$ws.Range("N51").Value = -19020

$ws.Range("H52").Value = 5000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 5000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 5000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -5452

$ws.Range("H81").Value = 917.53845
$ws.Range("I81").Value = 521.4286
$ws.Range("J81").Value = 1379.6666
$ws.Range("K81").Value = 1042.8572
$ws.Range("L81").Value = 2759.3332
$ws.Range("M81").Value = 18.14280000000008
$ws.Range("N81").Value = -4881.3332

$ws.Range("H84").Value = 917.53845
$ws.Range("I84").Value = 521.4286
$ws.Range("J84").Value = 1379.6666
$ws.Range("K84").Value = 5214.286
$ws.Range("L84").Value = 13796.666
$ws.Range("M84").Value = 89.71399999999994
$ws.Range("N84").Value = -24404.666

$ws.Range("H112").Value = 119729
$ws.Range("J112").Value = 119729
$ws.Range("L112").Value = 119729
$ws.Range("N112").Value = -122683

$ws.Range("H121").Value = 30420
$ws.Range("J121").Value = 30420
$ws.Range("L121").Value = 30420
$ws.Range("N121").Value = -33914
